$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New protocol entry (row 20)
$ws.Rows.Item(20).RowHeight = 64
$ws.Range("A20").Value = 45455
$ws.Range("A20").NumberFormat = "[$-F800]dddd\,\ mmmm\ dd\,\ yyyy"
$ws.Range("A20").HorizontalAlignment = -4131
$ws.Range("A20").VerticalAlignment = -4108
$ws.Range("A20").Font.Name = "Aptos Narrow"

$ws.Range("B20").Value = "Bis hierhin GRM neu digitalisiert (neue Gleichung, einfacher und representativer nach literatur)`nGPU code zur verarbeitung von rasterdaten geschrieben. Problem: Speicherbandbreite ist limitierender Faktor, nicht CPU`nGefahrenkarte mittels GRM"
$ws.Range("B20").WrapText = $true
$ws.Range("B20").VerticalAlignment = -4160

$ws.Range("D20").Select() | Out-Null
